# Applies the cryptos list refresh described by the commit:
# "Updated cryptos list on Thu Mar 28 07:06:40 UTC 2024 with GitHub Actions"
#
# D (Price) and E (Volume(1h)) columns hold plain-text values (no numeric
# NumberFormat in the workbook), including values that LOOK like numbers
# (e.g. "1.00", "0.998"). Writing such a string straight into .Value lets
# Excel auto-coerce it into a real number (dropping formatting / trailing
# zeros and tagging the cell s="2"+numFmtId=49). To avoid that we temporarily
# force the cell to Text ("@") before the write, then restore the default
# "Normal" style so the cell format matches the original (no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "70.339.27"; ForceText = $true },
    @{ Cell = "E2"; Value = "  +0.06%  "; ForceText = $false },
    @{ Cell = "D3"; Value = "3.572.72"; ForceText = $true },
    @{ Cell = "E3"; Value = "  -0.88%  "; ForceText = $false },
    @{ Cell = "E4"; Value = "  -0.09%  "; ForceText = $false },
    @{ Cell = "D5"; Value = "588.41"; ForceText = $true },
    @{ Cell = "E5"; Value = "  +1.70%  "; ForceText = $false },
    @{ Cell = "D6"; Value = "186.45"; ForceText = $true },
    @{ Cell = "E6"; Value = "  -1.62%  "; ForceText = $false },
    @{ Cell = "D7"; Value = "3.561.14"; ForceText = $true },
    @{ Cell = "E7"; Value = "  -1.08%  "; ForceText = $false },
    @{ Cell = "D8"; Value = "0.620"; ForceText = $true },
    @{ Cell = "E8"; Value = "  -1.37%  "; ForceText = $false },
    @{ Cell = "E9"; Value = "  +0.02%  "; ForceText = $false },
    @{ Cell = "D10"; Value = "0.201"; ForceText = $true },
    @{ Cell = "E10"; Value = "  +6.53%  "; ForceText = $false },
    @{ Cell = "E11"; Value = "  -1.48%  "; ForceText = $false },
    @{ Cell = "D12"; Value = "54.57"; ForceText = $true },
    @{ Cell = "E12"; Value = "  -2.50%  "; ForceText = $false },
    @{ Cell = "D13"; Value = "0.0000308"; ForceText = $true },
    @{ Cell = "E13"; Value = "  -1.90%  "; ForceText = $false },
    @{ Cell = "D14"; Value = "9.51"; ForceText = $true },
    @{ Cell = "E14"; Value = "  -1.54%  "; ForceText = $false },
    @{ Cell = "D15"; Value = "4.139.42"; ForceText = $true },
    @{ Cell = "E15"; Value = "  -1.14%  "; ForceText = $false },
    @{ Cell = "D16"; Value = "19.49"; ForceText = $true },
    @{ Cell = "E16"; Value = "  -1.63%  "; ForceText = $false },
    @{ Cell = "D17"; Value = "70.302.31"; ForceText = $true },
    @{ Cell = "E17"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D18"; Value = "3.566.31"; ForceText = $true },
    @{ Cell = "E18"; Value = "  -1.07%  "; ForceText = $false },
    @{ Cell = "D19"; Value = "12.48"; ForceText = $true },
    @{ Cell = "E19"; Value = "  -1.59%  "; ForceText = $false },
    @{ Cell = "E20"; Value = "  -1.10%  "; ForceText = $false },
    @{ Cell = "D21"; Value = "554.93"; ForceText = $true },
    @{ Cell = "E21"; Value = "  +12.43%  "; ForceText = $false },
    @{ Cell = "E22"; Value = "  -2.09%  "; ForceText = $false },
    @{ Cell = "D23"; Value = "17.88"; ForceText = $true },
    @{ Cell = "E23"; Value = "  -7.97%  "; ForceText = $false },
    @{ Cell = "D24"; Value = "4.67"; ForceText = $true },
    @{ Cell = "E24"; Value = "  +7.08%  "; ForceText = $false },
    @{ Cell = "E25"; Value = "  -0.78%  "; ForceText = $false },
    @{ Cell = "D26"; Value = "96.11"; ForceText = $true },
    @{ Cell = "E26"; Value = "  +0.08%  "; ForceText = $false },
    @{ Cell = "B27"; Value = "ImmutableX"; ForceText = $false },
    @{ Cell = "C27"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; ForceText = $false },
    @{ Cell = "D27"; Value = "3.00"; ForceText = $true },
    @{ Cell = "E27"; Value = "  +0.35%  "; ForceText = $false },
    @{ Cell = "B28"; Value = "RenderToken"; ForceText = $false },
    @{ Cell = "C28"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; ForceText = $false },
    @{ Cell = "D28"; Value = "11.21"; ForceText = $true },
    @{ Cell = "E28"; Value = "  +0.61%  "; ForceText = $false },
    @{ Cell = "E29"; Value = "  -2.84%  "; ForceText = $false },
    @{ Cell = "D30"; Value = "32.33"; ForceText = $true },
    @{ Cell = "E30"; Value = "  +0.62%  "; ForceText = $false },
    @{ Cell = "D31"; Value = "7.34"; ForceText = $true },
    @{ Cell = "E31"; Value = "  -3.93%  "; ForceText = $false },
    @{ Cell = "E32"; Value = "  +3.26%  "; ForceText = $false },
    @{ Cell = "D33"; Value = "65.20"; ForceText = $true },
    @{ Cell = "E33"; Value = "  -1.60%  "; ForceText = $false },
    @{ Cell = "E34"; Value = "  -2.72%  "; ForceText = $false },
    @{ Cell = "D35"; Value = "554.35"; ForceText = $true },
    @{ Cell = "E35"; Value = "  -3.48%  "; ForceText = $false },
    @{ Cell = "D36"; Value = "3.19"; ForceText = $true },
    @{ Cell = "E36"; Value = "  +4.33%  "; ForceText = $false },
    @{ Cell = "D37"; Value = "0.414"; ForceText = $true },
    @{ Cell = "E37"; Value = "  +4.59%  "; ForceText = $false },
    @{ Cell = "D38"; Value = "38.61"; ForceText = $true },
    @{ Cell = "E38"; Value = "  +0.10%  "; ForceText = $false },
    @{ Cell = "D39"; Value = "1.00"; ForceText = $true },
    @{ Cell = "E39"; Value = "  -0.08%  "; ForceText = $false },
    @{ Cell = "D40"; Value = "0.0₃0768"; ForceText = $true },
    @{ Cell = "E40"; Value = "  -5.51%  "; ForceText = $false },
    @{ Cell = "D41"; Value = "3.40"; ForceText = $true },
    @{ Cell = "E41"; Value = "  -4.43%  "; ForceText = $false },
    @{ Cell = "D42"; Value = "0.135"; ForceText = $true },
    @{ Cell = "E42"; Value = "  -2.05%  "; ForceText = $false },
    @{ Cell = "D43"; Value = "3.361.39"; ForceText = $true },
    @{ Cell = "E43"; Value = "  +3.89%  "; ForceText = $false },
    @{ Cell = "E44"; Value = "  -7.11%  "; ForceText = $false },
    @{ Cell = "D45"; Value = "3.54"; ForceText = $true },
    @{ Cell = "E45"; Value = "  +3.63%  "; ForceText = $false },
    @{ Cell = "E46"; Value = "  -1.18%  "; ForceText = $false },
    @{ Cell = "D47"; Value = "0.0444"; ForceText = $true },
    @{ Cell = "E47"; Value = "  +0.72%  "; ForceText = $false },
    @{ Cell = "D48"; Value = "9.22"; ForceText = $true },
    @{ Cell = "E48"; Value = "  -6.18%  "; ForceText = $false },
    @{ Cell = "D49"; Value = "0.136"; ForceText = $true },
    @{ Cell = "E49"; Value = "  -0.93%  "; ForceText = $false },
    @{ Cell = "D50"; Value = "0.998"; ForceText = $true },
    @{ Cell = "E50"; Value = "  -0.01%  "; ForceText = $false },
    @{ Cell = "D51"; Value = "137.47"; ForceText = $true },
    @{ Cell = "E51"; Value = "  +1.15%  "; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = "@"
        $rng.Value = $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}

